# "made expense sheet1 work"
# Replace the Walmart orders (rows 3 and 6) with new vendors/items so the
# expense sheet reflects the correct data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: Walmart -> Uber Eats (qty values stay the same: 1, 2, 1)
$ws.Range("A3").Value = "Uber Eats"
$ws.Range("D3").Value = "Burger"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "Red Wine"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = "Beef Steak"
$ws.Range("I3").Value = 1

# Row 6: Walmart -> Carvana
$ws.Range("A6").Value = "Carvana"
$ws.Range("D6").Value = "Truck"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "SUV"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "Car"
$ws.Range("I6").Value = 1

# Leave the selection where the author last clicked after editing.
$ws.Range("I7").Select()
